# lesson-06.pptx update:
#  1. Append a new "DOM & SAX" slide (Title + Content layout) at the end
#     of the deck.
#  2. Tidy up the eXtensible Markup Language title on slide 4 (drop the
#     stray trailing endParaRPr / code-style cleanup).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. New slide 20 - "DOM & SAX"
# ---------------------------------------------------------------------
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "DOM & SAX"

# Body / content placeholder
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "DOM – Document Object Model`rXmlDocument class`rSAX – Simple API for XML`rXmlReader"

# Second-level (indented) bullets
$body.Paragraphs(2).IndentLevel = 2
$body.Paragraphs(4).IndentLevel = 2

# Highlight the acronym letters in red, matching the rest of the deck
$p1 = $body.Paragraphs(1)
$p1.Characters(7, 1).Font.Color.RGB = 255
$p1.Characters(16, 1).Font.Color.RGB = 255
$p1.Characters(23, 1).Font.Color.RGB = 255

$p3 = $body.Paragraphs(3)
$p3.Characters(7, 1).Font.Color.RGB = 255
$p3.Characters(14, 1).Font.Color.RGB = 255
$p3.Characters(22, 1).Font.Color.RGB = 255

# ---------------------------------------------------------------------
# 2. Slide 4 - drop the stray trailing endParaRPr on the title
# ---------------------------------------------------------------------
$titleSlide4 = $p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange
$titleSlide4.Delete()
$titleSlide4.Text = "eXtensible Markup Language"
$titleSlide4.Characters(2, 1).Font.Color.RGB = 255
$titleSlide4.Characters(12, 1).Font.Color.RGB = 255
$titleSlide4.Characters(19, 1).Font.Color.RGB = 255
